# Generate Report for Handback
# Refresh the "Latest HO / Handback" timestamps for the second tracked file
# (c0848784-12d9-43de-9a7a-2d97d0b82dd1.md) across the Overview, zh-cn and
# de-de sheets after a new handback round completed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to c0848784-12d9-43de-9a7a-2d97d0b82dd1.md
$wsOverview.Range("G3").Value = "2016-09-05 23:01:22"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to c0848784-12d9-43de-9a7a-2d97d0b82dd1.md
$wsZhCn.Range("H3").Value = "2016-09-05 23:01:18"
$wsZhCn.Range("K3").Value = "2016-09-05 23:01:34"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to c0848784-12d9-43de-9a7a-2d97d0b82dd1.md
$wsDeDe.Range("H3").Value = "2016-09-05 23:01:22"
$wsDeDe.Range("K3").Value = "2016-09-05 23:01:42"
